# Auto-generated edit script: update Leve profit calculations across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 386
$ws.Range("I2").Value = 140
$ws.Range("K2").Value = 140
$ws.Range("M2").Value = -27
$ws.Range("H41").Value = 605.1429000000001
$ws.Range("J41").Value = 608.8
$ws.Range("L41").Value = 608.8
$ws.Range("N41").Value = -1488.8
$ws.Range("H58").Value = 1956
$ws.Range("I58").Value = 1956
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 5868
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -5718
$ws.Range("N58").ClearContents()
$ws.Range("H87").Value = 68316.664
$ws.Range("J87").Value = 68316.664
$ws.Range("L87").Value = 68316.664
$ws.Range("N87").Value = -70812.664
$ws.Range("H90").Value = 68316.664
$ws.Range("J90").Value = 68316.664
$ws.Range("L90").Value = 204949.992
$ws.Range("N90").Value = -217429.992
$ws.Range("H98").Value = 2207.3572
$ws.Range("I98").Value = 801.25
$ws.Range("K98").Value = 801.25
$ws.Range("M98").Value = 696.75
$ws.Range("H103").Value = 1936.6154
$ws.Range("I103").Value = 1539
$ws.Range("J103").Value = 2008.909
$ws.Range("K103").Value = 4617
$ws.Range("L103").Value = 6026.727000000001
$ws.Range("M103").Value = -4031
$ws.Range("N103").Value = -7198.727000000001
$ws.Range("H122").Value = 2207.3572
$ws.Range("I122").Value = 801.25
$ws.Range("K122").Value = 2403.75
$ws.Range("M122").Value = 46.25
$ws.Range("H132").Value = 1993.5834
$ws.Range("I132").Value = 1993.5834
$ws.Range("K132").Value = 5980.7502
$ws.Range("M132").Value = -3450.7502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8200.5
$ws.Range("I32").Value = 8200.5
$ws.Range("K32").Value = 8200.5
$ws.Range("M32").Value = -7913.5
$ws.Range("H88").Value = 2522.111
$ws.Range("J88").Value = 2980
$ws.Range("L88").Value = 2980
$ws.Range("N88").Value = -3792
$ws.Range("H91").Value = 2522.111
$ws.Range("J91").Value = 2980
$ws.Range("L91").Value = 2980
$ws.Range("N91").Value = -5788
$ws.Range("H102").Value = 2757.6667
$ws.Range("I102").Value = 2852.375
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2852.375
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -1230.375
$ws.Range("N102").Value = -5244
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 54701.668
$ws.Range("J100").Value = 54701.668
$ws.Range("L100").Value = 54701.668
$ws.Range("N100").Value = -56865.668
$ws.Range("H103").Value = 93999
$ws.Range("J103").Value = 93999
$ws.Range("L103").Value = 93999
$ws.Range("N103").Value = -96343
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2495.1667
$ws.Range("I31").Value = 1699.6666
$ws.Range("J31").Value = 3290.6667
$ws.Range("K31").Value = 1699.6666
$ws.Range("L31").Value = 3290.6667
$ws.Range("M31").Value = -1404.6666
$ws.Range("N31").Value = -3880.6667
$ws.Range("H34").Value = 2495.1667
$ws.Range("I34").Value = 1699.6666
$ws.Range("J34").Value = 3290.6667
$ws.Range("K34").Value = 1699.6666
$ws.Range("L34").Value = 3290.6667
$ws.Range("M34").Value = -1497.6666
$ws.Range("N34").Value = -3694.6667
$ws.Range("H134").Value = 2379.611
$ws.Range("I134").Value = 2011.5454
$ws.Range("K134").Value = 6034.6362
$ws.Range("M134").Value = -3499.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 706.875
$ws.Range("I12").Value = 1505
$ws.Range("J12").Value = 440.83334
$ws.Range("K12").Value = 4515
$ws.Range("L12").Value = 1322.50002
$ws.Range("M12").Value = -4342
$ws.Range("N12").Value = -1668.50002
$ws.Range("H97").Value = 1219.2858
$ws.Range("I97").Value = 1326.7273
$ws.Range("J97").Value = 825.3333
$ws.Range("K97").Value = 3980.1819
$ws.Range("L97").Value = 2475.9999
$ws.Range("M97").Value = -3484.1819
$ws.Range("N97").Value = -3467.9999
$ws.Range("H107").Value = 587.5833
$ws.Range("I107").Value = 450
$ws.Range("J107").Value = 615.1
$ws.Range("K107").Value = 1350
$ws.Range("L107").Value = 1845.3
$ws.Range("M107").Value = 570
$ws.Range("N107").Value = -5685.3
$ws.Range("H122").Value = 1793.7778
$ws.Range("I122").Value = 969.25
$ws.Range("J122").Value = 2453.4
$ws.Range("K122").Value = 8723.25
$ws.Range("L122").Value = 22080.6
$ws.Range("M122").Value = -6273.25
$ws.Range("N122").Value = -26980.6
$ws.Range("H128").Value = 390000
$ws.Range("I128").Value = 390000
$ws.Range("K128").Value = 1170000
$ws.Range("M128").Value = -1165020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3680.3572
$ws.Range("I132").Value = 2441.1428
$ws.Range("J132").Value = 4919.5713
$ws.Range("K132").Value = 7323.428400000001
$ws.Range("L132").Value = 14758.7139
$ws.Range("M132").Value = -4793.428400000001
$ws.Range("N132").Value = -19818.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4723.6665
$ws.Range("I81").Value = 2585.5
$ws.Range("K81").Value = 5171
$ws.Range("M81").Value = -4110
$ws.Range("H84").Value = 4723.6665
$ws.Range("I84").Value = 2585.5
$ws.Range("K84").Value = 25855
$ws.Range("M84").Value = -20551
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 1039.7273
$ws.Range("I100").Value = 1183.3334
$ws.Range("J100").Value = 867.4
$ws.Range("K100").Value = 2366.6668
$ws.Range("L100").Value = 1734.8
$ws.Range("M100").Value = -1825.6668
$ws.Range("N100").Value = -2816.8
$ws.Range("H113").Value = 1147.7
$ws.Range("I113").Value = 1599.6
$ws.Range("K113").Value = 4798.799999999999
$ws.Range("M113").Value = -2628.799999999999
$ws.Range("H122").Value = 1185.2
$ws.Range("I122").Value = 1181.5
$ws.Range("K122").Value = 3544.5
$ws.Range("M122").Value = -1094.5
$ws.Range("H136").Value = 888.46155
$ws.Range("I136").Value = 888.5
$ws.Range("K136").Value = 2665.5
$ws.Range("M136").Value = -115.5
